$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update only the changed odds columns ---
$ws.Range("G2").Value = 2.15
$ws.Range("I2").Value = 3.8
$ws.Range("J2").Value = 2.88
$ws.Range("L2").Value = 4.33
$ws.Range("M2").Value = 1.1
$ws.Range("N2").Value = 7
$ws.Range("U2").Value = 2.05
$ws.Range("V2").Value = 1.7
$ws.Range("W2").Value = 6
$ws.Range("X2").Value = 9
$ws.Range("Y2").Value = 9.5
$ws.Range("Z2").Value = 19
$ws.Range("AF2").Value = 67
$ws.Range("AG2").Value = 501
$ws.Range("AH2").Value = 9
$ws.Range("AI2").Value = 17
$ws.Range("AJ2").Value = 13
$ws.Range("AK2").Value = 41
$ws.Range("AL2").Value = 34
$ws.Range("AN2").Value = 4
$ws.Range("AO2").Value = 13
$ws.Range("AQ2").Value = 41
$ws.Range("AU2").Value = 9
$ws.Range("AW2").Value = 5.5
$ws.Range("AX2").Value = 21
$ws.Range("AY2").Value = 34
$ws.Range("AZ2").Value = 81
$ws.Range("BA2").Value = 126
$ws.Range("BB2").Value = 301

# --- Row 3: replaced entirely with new match data (previously row 4 content) ---
$ws.Range("A3").Value = "ClH8Rpbr"
$ws.Range("B3").Value = "27/11/2024"
$ws.Range("C3").Value = "13:00"
$ws.Range("D3").Value = "GEORGIA - CRYSTALBET EROVNULI LIGA"
$ws.Range("E3").Value = "Iberia 1999"
$ws.Range("F3").Value = "Kolkheti 1913"
$ws.Range("G3").Value = 1.22
$ws.Range("H3").Value = 5.4
$ws.Range("I3").Value = 9.25
$ws.Range("J3").Value = 1.6
$ws.Range("K3").Value = 2.72
$ws.Range("L3").Value = 7.4
$ws.Range("M3").Value = 1.01
$ws.Range("N3").Value = 17.5
$ws.Range("O3").Value = 1.08
$ws.Range("P3").Value = 6.5
$ws.Range("Q3").Value = 1.6
$ws.Range("R3").Value = 2.07
$ws.Range("S3").Value = 1.18
$ws.Range("T3").Value = 4.3
$ws.Range("U3").Value = 1.7
$ws.Range("V3").Value = 2.1
$ws.Range("W3").Value = 9
$ws.Range("X3").Value = 6.8
$ws.Range("Y3").Value = 7.8
$ws.Range("Z3").Value = 7.2
$ws.Range("AA3").Value = 8.25
$ws.Range("AB3").Value = 17.5
$ws.Range("AC3").Value = 20
$ws.Range("AD3").Value = 10.5
$ws.Range("AE3").Value = 16.5
$ws.Range("AF3").Value = 50
$ws.Range("AG3").Value = 250
$ws.Range("AH3").Value = 29
$ws.Range("AI3").Value = 65
$ws.Range("AJ3").Value = 24
$ws.Range("AK3").Value = 200
$ws.Range("AL3").Value = 80
$ws.Range("AM3").Value = 55
$ws.Range("AN3").Value = 3.35
$ws.Range("AO3").Value = 5.2
$ws.Range("AP3").Value = 12.5
$ws.Range("AQ3").Value = 11.75
$ws.Range("AR3").Value = 28
$ws.Range("AS3").Value = 120
$ws.Range("AT3").Value = 3.95
$ws.Range("AU3").Value = 8
$ws.Range("AV3").Value = 55
$ws.Range("AW3").Value = 10.5
$ws.Range("AX3").Value = 50
$ws.Range("AY3").Value = 40
$ws.Range("AZ3").Value = 350
$ws.Range("BA3").Value = 300
$ws.Range("BB3").Value = 400

# --- Row 4: replaced entirely with new match data ---
$ws.Range("A4").Value = "bsyYXnLE"
$ws.Range("B4").Value = "27/11/2024"
$ws.Range("C4").Value = "13:00"
$ws.Range("D4").Value = "GEORGIA - CRYSTALBET EROVNULI LIGA"
$ws.Range("E4").Value = "Torpedo Kutaisi"
$ws.Range("F4").Value = "Gagra"
$ws.Range("G4").Value = 1.44
$ws.Range("H4").Value = 3.7
$ws.Range("I4").Value = 6.9
$ws.Range("J4").Value = 1.93
$ws.Range("K4").Value = 2.18
$ws.Range("L4").Value = 6.4
$ws.Range("M4").Value = 1.01
$ws.Range("N4").Value = 11
$ws.Range("O4").Value = 1.15
$ws.Range("P4").Value = 4.05
$ws.Range("Q4").Value = 1.65
$ws.Range("R4").Value = 2
$ws.Range("S4").Value = 1.37
$ws.Range("T4").Value = 2.5
$ws.Range("U4").Value = 1.81
$ws.Range("V4").Value = 1.95
$ws.Range("W4").Value = 6.3
$ws.Range("X4").Value = 6.3
$ws.Range("Y4").Value = 6.6
$ws.Range("Z4").Value = 9
$ws.Range("AA4").Value = 9.25
$ws.Range("AB4").Value = 17.5
$ws.Range("AC4").Value = 11
$ws.Range("AD4").Value = 6.5
$ws.Range("AE4").Value = 12.5
$ws.Range("AF4").Value = 50
$ws.Range("AG4").Value = 300
$ws.Range("AH4").Value = 16
$ws.Range("AI4").Value = 40
$ws.Range("AJ4").Value = 17
$ws.Range("AK4").Value = 120
$ws.Range("AL4").Value = 60
$ws.Range("AM4").Value = 45
$ws.Range("AN4").Value = 3.3
$ws.Range("AO4").Value = 6.8
$ws.Range("AP4").Value = 14.5
$ws.Range("AQ4").Value = 20
$ws.Range("AR4").Value = 45
$ws.Range("AS4").Value = 175
$ws.Range("AT4").Value = 2.67
$ws.Range("AU4").Value = 7.2
$ws.Range("AV4").Value = 60
$ws.Range("AW4").Value = 8.25
$ws.Range("AX4").Value = 40
$ws.Range("AY4").Value = 37
$ws.Range("AZ4").Value = 300
$ws.Range("BA4").Value = 300
$ws.Range("BB4").Value = 500

# --- Remove old row 5 (data consolidated into rows 2-4; dimension becomes A1:BD4) ---
$ws.Rows("5:5").Delete()
